# Update cryptocurrency snapshot values (Price and Volume(1h) columns)
# as captured by the GitHub Actions scraper run on Sat Jan 14 19:33:23 UTC 2023.
# Cell values are plain text (e.g. "303.66", "5.61%") in the source sheet, so we
# force a text number format before assigning them and then clear the formatting
# override again so no stray style survives in the saved cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "303.66"
    "E2" = "5.61%"
    "D3" = "32.07"
    "E3" = "10.20%"
    "D4" = "5.271"
    "E4" = "1.46%"
    "E5" = "7.42%"
    "D6" = "7.847"
    "E6" = "5.68%"
    "D7" = "3.810"
    "E7" = "7.24%"
    "D8" = "1.518"
    "E8" = "8.55%"
    "D9" = "0.9174"
    "E9" = "2.05%"
    "D10" = "0.1686"
    "E10" = "5.05%"
    "D11" = "0.08024"
    "E11" = "6.63%"
    "D12" = "0.08063"
    "E12" = "4.52%"
    "D13" = "0.03012"
    "E13" = "2.63%"
    "D14" = "0.09905"
    "E14" = "9.94%"
    "D15" = "0.001492"
    "E15" = "-5.52%"
    "D16" = "0.04597"
    "E16" = "1.51%"
    "D17" = "0.006375"
    "E17" = "-1.53%"
    "E18" = "-0.36%"
    "D19" = "2.231"
    "E19" = "0.05%"
    "D20" = "0.3302"
    "E20" = "1.75%"
    "D21" = "0.1343"
    "E21" = "0.47%"
    "D22" = "4.486"
    "E22" = "11.23%"
    "E23" = "1.30%"
    "D24" = "0.001217"
    "D25" = "0.004446"
    "E25" = "7.24%"
    "E26" = "19.58%"
    "E27" = "6.12%"
    "D39" = "0.01716"
    "E39" = "2,545.34%"
    "D40" = "0.04509"
    "E40" = "3.33%"
    "D41" = "0.007149"
    "E41" = "3.21%"
    "D42" = "0.1347"
    "E42" = "8.07%"
    "D43" = "0.002256"
    "E43" = "9.11%"
    "D44" = "0.01285"
    "E44" = "9.93%"
    "D45" = "0.00006173"
    "E45" = "6.42%"
    "D46" = "0.7093"
    "E46" = "-63.23%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
